$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1877.4681
$ws.Range("I40").Value = 1615.7368
$ws.Range("J40").Value = 2055.0715
$ws.Range("K40").Value = 1615.7368
$ws.Range("L40").Value = 2055.0715
$ws.Range("M40").Value = -1440.7368
$ws.Range("N40").Value = -2405.0715

$ws.Range("H43").Value = 1404.65
$ws.Range("I43").Value = 2623.875
$ws.Range("J43").Value = 591.8333
$ws.Range("K43").Value = 2623.875
$ws.Range("L43").Value = 591.8333
$ws.Range("M43").Value = -2554.875
$ws.Range("N43").Value = -729.8333

$ws.Range("H51").Value = 2936.2727
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 3287.375
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 3287.375
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -4255.375

$ws.Range("H99").Value = 380
$ws.Range("I99").Value = 275.2
$ws.Range("J99").Value = 484.8
$ws.Range("K99").Value = 825.5999999999999
$ws.Range("L99").Value = 1454.4
$ws.Range("M99").Value = 672.4000000000001
$ws.Range("N99").Value = -4450.4

$ws.Range("H112").Value = 1346.3077
$ws.Range("J112").Value = 1476.579
$ws.Range("L112").Value = 4429.737
$ws.Range("N112").Value = -6645.737

$ws.Range("H113").Value = 2745.762
$ws.Range("I113").Value = 2863.6667
$ws.Range("J113").Value = 2451
$ws.Range("K113").Value = 2863.6667
$ws.Range("L113").Value = 2451
$ws.Range("M113").Value = 390.3332999999998
$ws.Range("N113").Value = -8959

$ws.Range("H115").Value = 400
$ws.Range("I115").Value = 400
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1200
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H116").Value = 5189.364
$ws.Range("I116").Value = 5322.125
$ws.Range("J116").Value = 4835.3335
$ws.Range("K116").Value = 5322.125
$ws.Range("L116").Value = 4835.3335
$ws.Range("M116").Value = -1880.125
$ws.Range("N116").Value = -11719.3335

$ws.Range("H137").Value = 976.05
$ws.Range("I137").Value = 934.5333000000001
$ws.Range("J137").Value = 1100.6
$ws.Range("K137").Value = 2803.5999
$ws.Range("L137").Value = 3301.8
$ws.Range("M137").Value = -253.5999000000002
$ws.Range("N137").Value = -8401.799999999999

$ws.Range("H138").Value = 4279.98
$ws.Range("I138").Value = 1042.2354
$ws.Range("J138").Value = 4951.2197
$ws.Range("K138").Value = 3126.7062
$ws.Range("L138").Value = 14853.6591
$ws.Range("M138").Value = 2013.2938
$ws.Range("N138").Value = -25133.6591

$ws.Range("H141").Value = 3411.8667
$ws.Range("I141").Value = 1115.7273
$ws.Range("J141").Value = 9726.25
$ws.Range("K141").Value = 3347.1819
$ws.Range("L141").Value = 29178.75
$ws.Range("M141").Value = 1832.8181
$ws.Range("N141").Value = -39538.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4236974.5
$ws.Range("I32").Value = 4850022.5
$ws.Range("J32").Value = 6941.4
$ws.Range("K32").Value = 4850022.5
$ws.Range("L32").Value = 6941.4
$ws.Range("M32").Value = -4849735.5
$ws.Range("N32").Value = -7515.4

$ws.Range("H132").Value = 2088
$ws.Range("I132").Value = 1054
$ws.Range("K132").Value = 3162
$ws.Range("M132").Value = -632

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 32750
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 32750
$ws.Range("K75").Value = 0
$ws.Range("N75").Value = -34622
$ws.Range("M75").ClearContents()

$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35630

$ws.Range("H78").Value = 32750
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 32750
$ws.Range("K78").Value = 0
$ws.Range("N78").Value = -107610
$ws.Range("M78").ClearContents()

$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37184

$ws.Range("H82").Value = 15524.5
$ws.Range("I82").Value = 11000
$ws.Range("J82").Value = 16655.625
$ws.Range("K82").Value = 11000
$ws.Range("L82").Value = 16655.625
$ws.Range("M82").Value = -10617
$ws.Range("N82").Value = -17421.625

$ws.Range("H85").Value = 15524.5
$ws.Range("I85").Value = 11000
$ws.Range("J85").Value = 16655.625
$ws.Range("K85").Value = 11000
$ws.Range("L85").Value = 16655.625
$ws.Range("M85").Value = -9674
$ws.Range("N85").Value = -19307.625

$ws.Range("H86").Value = 2070.8386
$ws.Range("I86").Value = 1891.5
$ws.Range("J86").Value = 2685.7144
$ws.Range("K86").Value = 1891.5
$ws.Range("L86").Value = 2685.7144
$ws.Range("M86").Value = -768.5
$ws.Range("N86").Value = -4931.7144

$ws.Range("H88").Value = 30000
$ws.Range("J88").Value = 30000
$ws.Range("L88").Value = 30000
$ws.Range("N88").Value = -30812

$ws.Range("H89").Value = 2070.8386
$ws.Range("I89").Value = 1891.5
$ws.Range("J89").Value = 2685.7144
$ws.Range("K89").Value = 9457.5
$ws.Range("L89").Value = 13428.572
$ws.Range("M89").Value = -3841.5
$ws.Range("N89").Value = -24660.572

$ws.Range("H91").Value = 30000
$ws.Range("J91").Value = 30000
$ws.Range("L91").Value = 30000
$ws.Range("N91").Value = -32808

$ws.Range("H99").Value = 43479892
$ws.Range("I99").Value = 62501730
$ws.Range("K99").Value = 62501730
$ws.Range("M99").Value = -62500232

$ws.Range("H134").Value = 46756.914
$ws.Range("I134").Value = 3517.6
$ws.Range("J134").Value = 80017.92
$ws.Range("K134").Value = 10552.8
$ws.Range("L134").Value = 240053.76
$ws.Range("M134").Value = -8017.799999999999
$ws.Range("N134").Value = -245123.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25496

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26716

$ws.Range("H94").Value = 2463.8518
$ws.Range("I94").Value = 2476.4285
$ws.Range("J94").Value = 2459.45
$ws.Range("K94").Value = 2476.4285
$ws.Range("L94").Value = 2459.45
$ws.Range("M94").Value = -2025.4285
$ws.Range("N94").Value = -3361.45

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 113.8
$ws.Range("I8").Value = 113.8
$ws.Range("K8").Value = 341.4
$ws.Range("M8").Value = -202.4

$ws.Range("H23").Value = 84.42856999999999
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 84.42856999999999
$ws.Range("K23").Value = 0
$ws.Range("N23").Value = -723.28571
$ws.Range("M23").ClearContents()

$ws.Range("H24").Value = 50
$ws.Range("I24").Value = 50
$ws.Range("K24").Value = 150
$ws.Range("M24").Value = 80

$ws.Range("H34").Value = 1750.6666
$ws.Range("J34").Value = 4986.6665
$ws.Range("L34").Value = 14959.9995
$ws.Range("N34").Value = -15127.9995

$ws.Range("H38").Value = 337.85715
$ws.Range("I38").Value = 58.75
$ws.Range("J38").Value = 710
$ws.Range("K38").Value = 176.25
$ws.Range("L38").Value = 2130
$ws.Range("M38").Value = 170.75
$ws.Range("N38").Value = -2824

$ws.Range("H46").Value = 2161
$ws.Range("J46").Value = 2650
$ws.Range("L46").Value = 7950
$ws.Range("N46").Value = -8132

$ws.Range("H55").Value = 25186.342
$ws.Range("I55").Value = 125485.375
$ws.Range("J55").Value = 871.42426
$ws.Range("K55").Value = 376456.125
$ws.Range("L55").Value = 2614.27278
$ws.Range("M55").Value = -376279.125
$ws.Range("N55").Value = -2968.27278

$ws.Range("H121").Value = 709
$ws.Range("J121").Value = 1007.3333
$ws.Range("L121").Value = 3021.9999
$ws.Range("N121").Value = -5641.9999

$ws.Range("H131").Value = 542930.8
$ws.Range("I131").Value = 8500412
$ws.Range("J131").Value = 12432.089
$ws.Range("K131").Value = 25501236
$ws.Range("L131").Value = 37296.267
$ws.Range("M131").Value = -25496196
$ws.Range("N131").Value = -47376.267

$ws.Range("H137").Value = 48385.273
$ws.Range("I137").Value = 2210.8333
$ws.Range("J137").Value = 103794.6
$ws.Range("K137").Value = 6632.499899999999
$ws.Range("L137").Value = 311383.8
$ws.Range("M137").Value = -1532.499899999999
$ws.Range("N137").Value = -321583.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4357.7856
$ws.Range("I70").Value = 3400
$ws.Range("J70").Value = 4889.8887
$ws.Range("K70").Value = 3400
$ws.Range("L70").Value = 4889.8887
$ws.Range("M70").Value = -3130
$ws.Range("N70").Value = -5429.8887

$ws.Range("H73").Value = 4357.7856
$ws.Range("I73").Value = 3400
$ws.Range("J73").Value = 4889.8887
$ws.Range("K73").Value = 3400
$ws.Range("L73").Value = 4889.8887
$ws.Range("M73").Value = -2464
$ws.Range("N73").Value = -6761.8887

$ws.Range("H113").Value = 4160.8667
$ws.Range("I113").Value = 5190.909
$ws.Range("J113").Value = 1328.25
$ws.Range("K113").Value = 5190.909
$ws.Range("L113").Value = 1328.25
$ws.Range("M113").Value = -3020.909
$ws.Range("N113").Value = -5668.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1204.2858
$ws.Range("I46").Value = 1168.4615
$ws.Range("J46").Value = 1262.5
$ws.Range("K46").Value = 1168.4615
$ws.Range("L46").Value = 1262.5
$ws.Range("M46").Value = -980.4614999999999
$ws.Range("N46").Value = -1638.5

$ws.Range("H132").Value = 1839.356
$ws.Range("I132").Value = 1674.6666
$ws.Range("J132").Value = 2048.3845
$ws.Range("K132").Value = 5023.9998
$ws.Range("L132").Value = 6145.1535
$ws.Range("M132").Value = -2493.9998
$ws.Range("N132").Value = -11205.1535

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1544
$ws.Range("I132").Value = 1052.3
$ws.Range("J132").Value = 2158.625
$ws.Range("K132").Value = 3156.9
$ws.Range("L132").Value = 6475.875
$ws.Range("M132").Value = -626.8999999999996
$ws.Range("N132").Value = -11535.875

$ws.Range("H136").Value = 1548.8784
$ws.Range("I136").Value = 1478.5574
$ws.Range("J136").Value = 1878.8462
$ws.Range("K136").Value = 4435.6722
$ws.Range("L136").Value = 5636.5386
$ws.Range("M136").Value = -1885.6722
$ws.Range("N136").Value = -10736.5386
